# Apply the "vegi_coop" crawl refresh:
#  - every row's timestamp (column O) moves from the 12:58:15 crawl to the
#    20:51:33 crawl
#  - a handful of ratingAmount values (column D) ticked up by one new review
#  - one product's aria-label (column M) picked up an "Online kein Bestand"
#    stock-out note

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2023-01-06 20:51:33"

$lastRow = 520
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 15).Value = $newTimestamp
}

# ratingAmount bumps (column D)
$ws.Range("D136").Value = 15
$ws.Range("D191").Value = 10
$ws.Range("D193").Value = 10
$ws.Range("D238").Value = 16
$ws.Range("D240").Value = 2

# productAriaLabel update for Betty Bossi Kürbis & Kichererbsen (row 519)
$ws.Range("M519").Value = "Betty Bossi Kürbis &amp; Kichererbsen - Online kein Bestand 5.50 Schweizer Franken"
